$d = $word.ActiveDocument

# 1) Replace the "Store the tables in S3 using encryption" bullet text
$d.Content.Find.Execute(
    "Store the tables in S3 using encryption",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "There is nothing to store that would require S3. ",
    2
)

# 2) Add a new bullet "Encrypt the data in DynamoDB" right after the
#    "Encrypt the data in S3" bullet, matching its list level/formatting.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Encrypt the data in S3") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.InsertParagraphAfter()
    $newRange = $target.Next().Range
    $newRange.Text = "Encrypt the data in DynamoDB"
}
